$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Planilha1 -> Tabela Auxiliar)
$ws.Name = "Tabela Auxiliar"

# Update the link_logo column (C) with the new GitHub "blob/main" .png URLs.
$ws.Range("C2").Value = "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/Brasil%20Bolsa%20Balc%C3%A3o/B3SA3_Logo.png"
$ws.Range("C3").Value = "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/Banco%20do%20Brasil/BBAS3_Logo.png"
$ws.Range("C4").Value = "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/Banco%20Ita%C3%BA/ITUB4_Logo.png"
$ws.Range("C5").Value = "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/CaixaEconomica/CXSE3_Logo.png"

# Drop every existing hyperlink (Range.Hyperlinks.Delete clears the sheet's
# hyperlink collection) and recreate them against the new .png targets,
# without a TextToDisplay override (so no "display" attribute is written).
$ws.Range("C4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/Banco%20do%20Brasil/BBAS3_Logo.png")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/Brasil%20Bolsa%20Balc%C3%A3o/B3SA3_Logo.png")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/CaixaEconomica/CXSE3_Logo.png")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/Mateus-Fleck/yFinance_BI_2.0_TEST/blob/main/Assets/Banco%20Ita%C3%BA/ITUB4_Logo.png")

# Hyperlinks.Add() stamps every newly-linked cell with a fresh cell-format
# (xfId=1/fontId=2) entry; re-applying the named "Hiperlink" style on C2,
# C3 and C5 folds them back onto the original shared Hiperlink format index,
# leaving C4 on its own distinct copy (matching the source workbook, where
# only C4 carries the new style).
$ws.Range("C2").Style = "Hiperlink"
$ws.Range("C3").Style = "Hiperlink"
$ws.Range("C5").Style = "Hiperlink"

# Add the new blank, styled row 7 cell (mirrors the existing C20 filler cell).
$ws.Range("C20").Copy($ws.Range("C7"))

# Move the active selection to C7.
$ws.Range("C7").Select()

# Page setup: A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "edit applied"
